$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cluster label values for column A (rows 2-26)
$newValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 3
    8  = 3
    9  = 3
    10 = 3
    11 = 3
    12 = 3
    13 = 3
    14 = 3
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 4
    20 = 4
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 1
    26 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("A$row").Value = $newValues[$row]
}
